$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two "Requisitos" entries: row 24 currently holds the
# "LOM3202 - Circuitos Elétricos (Indicação de Conjunto)" text and row 25
# holds the "LOB1053 - Física III (Requisito)" text. The new order puts
# the LOB1053 entry first (row 24) and the LOM3202 entry second (row 25).

$text1 = "LOM3202 -  Circuitos Elétricos  (Indicação de Conjunto)`n"
$text2 = "LOB1053 -  Física III  (Requisito)`n"

$ws.Range("B24").Value = $text2
$ws.Range("C24").Value = $text2

$ws.Range("B25").Value = $text1
$ws.Range("C25").Value = $text1
